$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.666.13'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('D3').Value = '''3.504.13'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''590.75'
$ws.Range('E5').Value = '  +3.27%  '
$ws.Range('D6').Value = '''170.16'
$ws.Range('E6').Value = '  +5.70%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '''3.500.47'
$ws.Range('E8').Value = '  +2.16%  '
$ws.Range('D9').Value = '''0.583'
$ws.Range('E9').Value = '  +4.54%  '
$ws.Range('D10').Value = '''7.35'
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').Value = '''0.126'
$ws.Range('E11').Value = '  +4.19%  '
$ws.Range('D12').Value = '''0.441'
$ws.Range('E12').Value = '  +3.55%  '
$ws.Range('D13').Value = '''4.106.01'
$ws.Range('E13').Value = '  +2.18%  '
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('D15').Value = '''28.39'
$ws.Range('E15').Value = '  +4.39%  '
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('D17').Value = '''66.648.86'
$ws.Range('E17').Value = '  +3.75%  '
$ws.Range('D18').Value = '''3.494.56'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').Value = '''6.35'
$ws.Range('E19').Value = '  +3.86%  '
$ws.Range('D20').Value = '''14.06'
$ws.Range('E20').Value = '  +3.14%  '
$ws.Range('D21').Value = '''391.13'
$ws.Range('E21').Value = '  +2.89%  '
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').Value = '''73.06'
$ws.Range('E23').Value = '  +2.12%  '
$ws.Range('D24').Value = '''1.00'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Value = '''0.535'
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('D26').Value = '''0.0000123'
$ws.Range('E26').Value = '  +5.39%  '
$ws.Range('D27').Value = '''10.30'
$ws.Range('E27').Value = '  +7.35%  '
$ws.Range('E28').Value = '  +2.58%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '''6.38'
$ws.Range('E30').Value = '  +5.00%  '
$ws.Range('D31').Value = '''1.49'
$ws.Range('E31').Value = '  +5.58%  '
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('D33').Value = '''23.61'
$ws.Range('E33').Value = '  +2.58%  '
$ws.Range('E34').Value = '  +5.12%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '''1.62'
$ws.Range('E36').Value = '  +6.62%  '
$ws.Range('D37').Value = '''162.69'
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('D38').Value = '''0.883'
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('D39').Value = '''1.91'
$ws.Range('E39').Value = '  +4.34%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''6.84'
$ws.Range('E40').Value = '  +5.60%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '''4.72'
$ws.Range('E41').Value = '  +5.50%  '
$ws.Range('D42').Value = '''0.0748'
$ws.Range('E42').Value = '  +2.08%  '
$ws.Range('D43').Value = '''26.59'
$ws.Range('E43').Value = '  +2.99%  '
$ws.Range('D44').Value = '''2.813.81'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').Value = '''26.85'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('D46').Value = '''43.02'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').Value = '''0.0313'
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '''2.52'
$ws.Range('E48').Value = '  +4.91%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '''354.07'
$ws.Range('E49').Value = '  +3.22%  '
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('D51').Value = '''33.79'
$ws.Range('E51').Value = '  +12.00%  '
